$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new attendance week (column K, date 2019-10-19) was recorded for the
# plus-education sign-in sheet. Every student attended that week except
# row 5 (오대완) and row 15 (최준우), who have no mark for that date.
$attendedRows = @(2,3,4,6,7,8,9,10,11,12,13,14,16,17,18,19,20)
foreach ($r in $attendedRows) {
    $ws.Cells.Item($r, 11).Value = 1
}

# The comment on L15 (지각 / late) gets an extra note appended about what
# time the student actually arrived.
$cmt = $ws.Range("L15").Comment
$cmt.Text("지각`n12시에 들어옴`n")

# Leave the selection where the editor was last working.
$ws.Range("L15").Select()
